$d = $word.ActiveDocument

# --- Locate the paragraph that owns the "_GoBack" bookmark. This is the
#     list item "Filtrare le lune per pianeta ... tramite click dalla
#     pagina di descrizione del pianeta" that must be merged into the
#     following "Dopo aver finito" paragraph, keeping only the bookmark.
$bm = $d.Bookmarks.Item("_GoBack")
$bmParaRange = $d.Range($bm.Start, $bm.Start)
$null = $bmParaRange.Expand(4)   # wdParagraph
$paraEnd = $bmParaRange.End

# --- Step 1: merge that paragraph with the next one ("Dopo aver finito")
#     by deleting the paragraph mark between them.
$markRange = $d.Range($paraEnd - 1, $paraEnd)
$markRange.Delete()

# --- Step 2: find "Dopo aver finito" text, now sharing the merged
#     paragraph, so we know where the text to keep begins.
$findRange = $d.Content
$null = $findRange.Find.Execute("Dopo aver finito", $false, $false, $false,
                                 $false, $false, $true, 1, $false, "", 0)
$keepStart = $findRange.Start

# --- Step 3: delete the leftover text after the bookmark (up to the text
#     we keep), then the leftover text before the bookmark (back to the
#     start of the paragraph). The bookmark itself (a collapsed range) is
#     left untouched, so it survives and ends up right before "Dopo aver
#     finito".
$bm = $d.Bookmarks.Item("_GoBack")
$suffix = $d.Range($bm.Start, $keepStart)
$suffix.Delete()

$bm = $d.Bookmarks.Item("_GoBack")
$paraStartRange = $d.Range($bm.Start, $bm.Start)
$null = $paraStartRange.Expand(4)   # wdParagraph
$paraStart = $paraStartRange.Start
$prefix = $d.Range($paraStart, $bm.Start)
$prefix.Delete()
